$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.680891275405884
$ws.Range("B1").Value = 2.100168228149414
$ws.Range("C1").Value = 5.103261947631836
$ws.Range("D1").Value = 1.397628664970398
$ws.Range("E1").Value = 0.6636487245559692
